# Apply updated values to the "Valores" worksheet to reflect only the
# "competencia" figures in the financeiro export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 4

# Row 9
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0

# Row 10
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 4
$ws.Range("L10").Value = 0

# Row 11
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0

# Row 13
$ws.Range("A13").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1

# Row 15
$ws.Range("A15").Value = 4
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 5
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 4

# Row 16
$ws.Range("A16").Value = 26
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 19
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 10
$ws.Range("L16").Value = 9
